$wb = $excel.ActiveWorkbook
$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux  = $wb.Worksheets.Item("Totaux")

# ---------------------------------------------------------------------------
# Journal sheet: turn the blank template row 13 into real data, and add a
# brand-new row 14 below it (copy row 13's formatting down first so the new
# row inherits the same cell styles used throughout the table).
# ---------------------------------------------------------------------------
$wsJournal.Range("A13:D13").Copy($wsJournal.Range("A14:D14"))

$wsJournal.Range("A13").Value = 44967
$wsJournal.Range("B13").Value = 2
$wsJournal.Range("C13").Value = 0.069444444444444434
$wsJournal.Range("D13").Value = "Analyse"
$wsJournal.Range("E13").Value = 'Analyse de la partie "gestion du parc" sur la nouvelle version de l''application '

$wsJournal.Range("A14").Value = 44967
$wsJournal.Range("B14").Value = 2
$wsJournal.Range("C14").Value = 0.03125
$wsJournal.Range("D14").Value = "Export"
$wsJournal.Range("E14").Value = "Export des données des catégories sur l'ancienne version de l'application"

# Grow "Tableau1" to include the new row.
$loJournal = $wsJournal.ListObjects.Item(1)
$loJournal.Resize($wsJournal.Range("A1:E14"))

# ---------------------------------------------------------------------------
# Totaux sheet: row 5's date shifts a day, the old "Total" row (row 6) turns
# into a regular weekly-sum data row, and a fresh "Total" row is appended as
# row 7.
# ---------------------------------------------------------------------------

# Preserve the "Total" row styling by copying it one row down first.
$wsTotaux.Range("A6:B6").Copy($wsTotaux.Range("A7:B7"))
# Then restyle row 6 like an ordinary date/value row (copy row 5's look).
$wsTotaux.Range("A5:B5").Copy($wsTotaux.Range("A6:B6"))

$wsTotaux.Range("A5").Value = 44966

$wsTotaux.Range("A6").Value = 44967
$wsTotaux.Range("B6").Formula = "=SUM(Journal!C13:C14)"

$wsTotaux.Range("B7").Formula = "=SUM(B2:B6)"

# Grow "Tableau2" to include the new row (totals row now sits on row 7).
$loTotaux = $wsTotaux.ListObjects.Item(1)
$loTotaux.Resize($wsTotaux.Range("A1:B7"))

# ---------------------------------------------------------------------------
# Selection bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$wsJournal.Activate() | Out-Null
$wsJournal.Range("E19").Select() | Out-Null

$wsTotaux.Activate() | Out-Null
$wsTotaux.Range("J23").Select() | Out-Null
